$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 34; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq 45641) {
        $cell.Value = 45642
    }
}
